{"js": "// Add four new customer rows to the \"Danh s\u00e1ch kh\u00e1ch h\u00e0ng\" list:\n//   1) a \"Nguyen Anh Thu\" row (CMND 1989) right before the 1st \"Nguy\u1ec5n V\u0103n L\u00fd\" row\n//   2) \"Nguy\u1ec5n V\u0103Ho\u00e0ng\" / \"Nguy\u1ec5n V\u0103n Ho\u00e0ng\" / \"Nguy\u1ec5n Duy Anh\" rows right before\n//      the 2nd (duplicate) \"Nguy\u1ec5n V\u0103n L\u00fd\" row\n//   3) a \"Nguy\u1ec5n V\u0103n T\u00e0i\" row right after the \"Nguyen Van Minh Sau Hai\" row\n//\n// Paragraphs are located by their (stable) text rather than a hard-coded index,\n// since that is resilient to any other edits already present in the body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Returns the paragraph object for the `occurrence`-th (0-based) paragraph\n// whose text contains `needle`.\nfunction findParagraphByText(needle, occurrence) {\n  let seen = 0;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.includes(needle)) {\n      if (seen === occurrence) return paragraphs.items[i];\n      seen++;\n    }\n  }\n  throw new Error(`Paragraph containing \"${needle}\" (occurrence ${occurrence}) not found`);\n}\n\n// --- Hunk 1: before the 1st \"Nguy\u1ec5n V\u0103n L\u00fd\" row -----------------------------\nconst ly1 = findParagraphByText(\"Nguy\u1ec5n V\u0103n L\u00fd\", 0);\nly1.insertParagraph(\n  \"H\u1ecd t\u00ean: Nguyen Anh Thu - S\u0110T: 0964561306 - \u0110\u1ecba ch\u1ec9: VN - Email: ndanhdev3@gmail.com - CMND: 1989 - N\u0103m sinh: 123456789\",\n  Word.InsertLocation.before\n);\n\n// --- Hunk 2: before the 2nd (duplicate) \"Nguy\u1ec5n V\u0103n L\u00fd\" row -----------------\nconst ly2 = findParagraphByText(\"Nguy\u1ec5n V\u0103n L\u00fd\", 1);\nly2.insertParagraph(\n  \"H\u1ecd t\u00ean: Nguy\u1ec5n V\u0103Ho\u00e0ng - S\u0110T: 0964561306 - \u0110\u1ecba ch\u1ec9: VN - Email: aduy644@gmail.com - CMND: 1989 - N\u0103m sinh: 123456789\",\n  Word.InsertLocation.before\n);\nly2.insertParagraph(\n  \"H\u1ecd t\u00ean: Nguy\u1ec5n V\u0103n Ho\u00e0ng - S\u0110T: 0964561306 - \u0110\u1ecba ch\u1ec9: VN - Email: aduy644@gmail.com - CMND: 1989 - N\u0103m sinh: 123456789\",\n  Word.InsertLocation.before\n);\nly2.insertParagraph(\n  \"H\u1ecd t\u00ean: Nguy\u1ec5n Duy Anh - S\u0110T: 0964561306 - \u0110\u1ecba ch\u1ec9: VN - Email: aduy644@gmail.com - CMND: 1989 - N\u0103m sinh: 123456789\",\n  Word.InsertLocation.before\n);\n\n// --- Hunk 3: after the \"Nguyen Van Minh Sau Hai\" row ------------------------\nconst sauHai = findParagraphByText(\"Nguyen Van Minh Sau Hai\", 0);\nsauHai.insertParagraph(\n  \"H\u1ecd t\u00ean: Nguy\u1ec5n V\u0103n T\u00e0i - S\u0110T: 0964561306 - \u0110\u1ecba ch\u1ec9: VN - Email: aduy644@gmail.com - CMND: 1989 - N\u0103m sinh: 123456789\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Add four new customer rows to the \"Danh s\u00e1ch kh\u00e1ch h\u00e0ng\" list:\n#   1) a \"Nguyen Anh Thu\" row (CMND 1989) right before the 1st \"Nguy\u1ec5n V\u0103n L\u00fd\" row\n#   2) \"Nguy\u1ec5n V\u0103Ho\u00e0ng\" / \"Nguy\u1ec5n V\u0103n Ho\u00e0ng\" / \"Nguy\u1ec5n Duy Anh\" rows right before\n#      the 2nd (duplicate) \"Nguy\u1ec5n V\u0103n L\u00fd\" row\n#   3) a \"Nguy\u1ec5n V\u0103n T\u00e0i\" row right after the \"Nguyen Van Minh Sau Hai\" row\n#\n# Paragraphs are located by their (stable) text rather than a hard-coded index,\n# since that is resilient to any other edits already present in the body.\n# New rows are always added with InsertParagraphAfter() on the paragraph that\n# should immediately precede them (rather than InsertParagraphBefore() on the\n# following paragraph) so the anchor paragraph's own Range is left untouched.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphIndexByText($needle, $occurrence) {\n    $seen = 0\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        if ($d.Paragraphs.Item($i).Range.Text -like \"*$needle*\") {\n            if ($seen -eq $occurrence) { return $i }\n            $seen = $seen + 1\n        }\n    }\n    throw \"Paragraph containing '$needle' (occurrence $occurrence) not found\"\n}\n\nfunction Add-RowAfter($paragraphIndex, $text) {\n    # Appends a new paragraph right after $paragraphIndex and returns the new\n    # paragraph's index.\n    $d.Paragraphs.Item($paragraphIndex).Range.InsertParagraphAfter()\n    $newIndex = $paragraphIndex + 1\n    $d.Paragraphs.Item($newIndex).Range.Text = $text\n    return $newIndex\n}\n\n# --- Hunk 1: before the 1st \"Nguy\u1ec5n V\u0103n L\u00fd\" row -----------------------------\n$ly1Index = Get-ParagraphIndexByText \"Nguy\u1ec5n V\u0103n L\u00fd\" 0\nAdd-RowAfter ($ly1Index - 1) \"H\u1ecd t\u00ean: Nguyen Anh Thu - S\u0110T: 0964561306 - \u0110\u1ecba ch\u1ec9: VN - Email: ndanhdev3@gmail.com - CMND: 1989 - N\u0103m sinh: 123456789\" | Out-Null\n\n# --- Hunk 2: before the 2nd (duplicate) \"Nguy\u1ec5n V\u0103n L\u00fd\" row -----------------\n$ly2Index = Get-ParagraphIndexByText \"Nguy\u1ec5n V\u0103n L\u00fd\" 1\n$lastNewIndex = Add-RowAfter ($ly2Index - 1) \"H\u1ecd t\u00ean: Nguy\u1ec5n V\u0103Ho\u00e0ng - S\u0110T: 0964561306 - \u0110\u1ecba ch\u1ec9: VN - Email: aduy644@gmail.com - CMND: 1989 - N\u0103m sinh: 123456789\"\n$lastNewIndex = Add-RowAfter $lastNewIndex \"H\u1ecd t\u00ean: Nguy\u1ec5n V\u0103n Ho\u00e0ng - S\u0110T: 0964561306 - \u0110\u1ecba ch\u1ec9: VN - Email: aduy644@gmail.com - CMND: 1989 - N\u0103m sinh: 123456789\"\n$lastNewIndex = Add-RowAfter $lastNewIndex \"H\u1ecd t\u00ean: Nguy\u1ec5n Duy Anh - S\u0110T: 0964561306 - \u0110\u1ecba ch\u1ec9: VN - Email: aduy644@gmail.com - CMND: 1989 - N\u0103m sinh: 123456789\"\n\n# --- Hunk 3: after the \"Nguyen Van Minh Sau Hai\" row ------------------------\n$sauHaiIndex = Get-ParagraphIndexByText \"Nguyen Van Minh Sau Hai\" 0\nAdd-RowAfter $sauHaiIndex \"H\u1ecd t\u00ean: Nguy\u1ec5n V\u0103n T\u00e0i - S\u0110T: 0964561306 - \u0110\u1ecba ch\u1ec9: VN - Email: aduy644@gmail.com - CMND: 1989 - N\u0103m sinh: 123456789\" | Out-Null\n\nWrite-Output \"done\"\n"}
